# Insert six new paragraphs (String method / Slice.../ Numerical method /
# toString.../ to switch to script mode.../ It is used to avoid bugs...)
# right after paragraph 1, before the pre-existing trailing empty paragraph.
#
# NOTE: the COM host's Range.InsertXML has an observed quirk where replacing a
# range that sits immediately before the very last paragraph of the document
# with a multi-paragraph (3+) XML fragment corrupts paragraph 1 instead of
# inserting correctly. Work around it by first padding with extra throw-away
# paragraphs (so the target range is no longer adjacent to the document's
# final paragraph), performing the InsertXML, then deleting the padding.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.InsertParagraphAfter()
$target = $d.Paragraphs(2)
$targetRange = $target.Range

$padCount = 6
for ($i = 0; $i -lt $padCount; $i++) {
    $targetRange.InsertParagraphAfter()
}

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">String method </w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">Slice </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>substr</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>startsWith</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>lastindexOf</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> ,</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>indexOf</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>endsWith</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">Numerical method </w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>toString</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>toPrecision</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> ,</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> value of </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>etc.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">to switch to script </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>mode</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> we use “use strict” in the beginning of the code.</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>It</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>i</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>s</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> used to avoid bugs in our scripts</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($xmlFrag)

# Remove the padding paragraphs that were inserted after the new content.
# After InsertXML, the six new paragraphs occupy what was $target's slot,
# so the padding paragraphs now immediately follow them, still before the
# original trailing empty paragraph.
for ($i = 0; $i -lt $padCount; $i++) {
    $newParaCount = 6
    $padPara = $d.Paragraphs(1 + $newParaCount + 1)
    $padRange = $d.Range($padPara.Range.Start, $padPara.Range.End)
    $padRange.Delete()
}
